# Aspects update: add a "TYPE" row/column to the Edges interaction matrix on
# the "Category usage" sheet (mirrors the "TYPE" row/column that already
# exists in the "Nodes" matrix), and add a new "c7" legend entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Category usage")

function Set-StyledValue {
    param(
        [string]$Dest,
        [string]$Value,
        [string]$DonorStyle
    )
    $ws.Range($Dest).Value2 = $Value
    if ($DonorStyle) {
        $ws.Range($DonorStyle).Copy()
        $ws.Range($Dest).PasteSpecial(-4122)
    }
}

function Set-StyledFormula {
    param(
        [string]$Dest,
        [string]$Formula,
        [string]$DonorStyle
    )
    $ws.Range($Dest).Formula = $Formula
    if ($DonorStyle) {
        $ws.Range($DonorStyle).Copy()
        $ws.Range($Dest).PasteSpecial(-4122)
    }
}

# --- New legend entry (c7) next to the PARAM block, row 35 -----------------
Set-StyledValue -Dest "S35" -Value "c7" -DonorStyle "S29"
Set-StyledValue -Dest "T35" -Value "c1 but for assignments, which are only for readers and creators" -DonorStyle "T29"

# --- Edges matrix header row (row 39): insert TYPE before MULT_IN ----------
Set-StyledValue -Dest "O39" -Value "TYPE" -DonorStyle "N39"
Set-StyledValue -Dest "P39" -Value "MULT_I-" -DonorStyle "N39"
Set-StyledValue -Dest "Q39" -Value "MULT_OUT" -DonorStyle "N39"
Set-StyledValue -Dest "R39" -Value "ASSOC" -DonorStyle "N39"

# --- REMARK row (row 40) ----------------------------------------------------
Set-StyledValue -Dest "O40" -Value "-" -DonorStyle "N40"
Set-StyledValue -Dest "R40" -Value "-" -DonorStyle "M40"

# --- SORT row (row 41) ------------------------------------------------------
Set-StyledValue -Dest "O41" -Value "-" -DonorStyle "N40"
Set-StyledValue -Dest "P41" -Value "-" -DonorStyle "M40"
Set-StyledValue -Dest "Q41" -Value "X" -DonorStyle "E30"
Set-StyledValue -Dest "R41" -Value "-" -DonorStyle "M40"

# --- LABEL row (row 42) -----------------------------------------------------
Set-StyledValue -Dest "Q42" -Value "X" -DonorStyle "E30"
Set-StyledValue -Dest "R42" -Value "X" -DonorStyle "B3"

# --- TYPE row (row 43, new row label; was MULT_IN) --------------------------
Set-StyledValue -Dest "K43" -Value "TYPE" -DonorStyle $null
Set-StyledValue -Dest "Q43" -Value "X" -DonorStyle "E30"
Set-StyledValue -Dest "R43" -Value "X" -DonorStyle "B3"

# --- MULT_IN row (row 44, was MULT_OUT) -------------------------------------
Set-StyledValue -Dest "K44" -Value "MULT_IN" -DonorStyle $null
Set-StyledValue -Dest "R44" -Value "c1" -DonorStyle "E30"

# --- MULT_OUT row (row 45, was ASSOC) ---------------------------------------
Set-StyledValue -Dest "K45" -Value "MULT_OUT" -DonorStyle $null
Set-StyledValue -Dest "R45" -Value "X" -DonorStyle "E30"

# --- ASSOC row (row 46, brand-new row) --------------------------------------
Set-StyledValue -Dest "K46" -Value "ASSOC" -DonorStyle $null
Set-StyledFormula -Dest "L46" -Formula "=R40" -DonorStyle "L45"
Set-StyledFormula -Dest "M46" -Formula "=R41" -DonorStyle "L45"
Set-StyledFormula -Dest "N46" -Formula "=R42" -DonorStyle "L45"
Set-StyledFormula -Dest "O46" -Formula "=R43" -DonorStyle "L45"
Set-StyledFormula -Dest "P46" -Formula "=R44" -DonorStyle "L45"
Set-StyledFormula -Dest "Q46" -Formula "=R45" -DonorStyle "L45"
$ws.Range("Q45").Copy()
$ws.Range("R46").PasteSpecial(-4122)

# --- Column width: extend the uniform 6.140625-wide column formatting from
# column Q (17) to include the new column R (18) ----------------------------
$ws.Range("R1").ColumnWidth = $ws.Range("Q1").ColumnWidth

# --- Update saved selection / scroll state to match the editing session ----
$ws.Range("R44").Select()
